# Actualizacion data y funcionalidades
# Updates the test-data values on the "Depositos" sheet and refreshes the
# active cell selection left by the author when the file was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Depositos")

# usuario: invictus10 -> pruebauser01
$ws.Range("D2").Value = "pruebauser01"

# numeroCuenta: 406-125210-01 -> 406-182800-02
# (leading apostrophe keeps the cell's existing "quote prefix" text style,
# matching the unchanged cell format in the target workbook)
$ws.Range("N2").Value = "'406-182800-02"

# Selection left on the sheet moved from M11 to M6
$ws.Range("M6").Select() | Out-Null
